# Fix tab names in input TMAs: rename "pten" sheet to "PTEN"
$wb = $excel.ActiveWorkbook

$tmaSheet  = $wb.Worksheets.Item("TMA map")
$ptenSheet = $wb.Worksheets.Item("pten")

$ptenSheet.Name = "PTEN"

# Restore/keep the selection on "TMA map" (F10), then leave it as a non-active tab
$tmaSheet.Activate() | Out-Null
$tmaSheet.Range("F10").Select() | Out-Null

# PTEN becomes the active/selected tab, with its selection moved to D3
$ptenSheet.Activate() | Out-Null
$ptenSheet.Range("D3").Select() | Out-Null
